$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns: lat / long
$ws.Range("G1").Value = "lat"
$ws.Range("H1").Value = "long"

# Row 2 - CC-407 (Eastern North America)
$ws.Range("G2").Value = 39.045752999999998
$ws.Range("H2").Value = -76.641272999999998

# Row 3 - CC-125 (Amherst, MA)
$ws.Range("G3").Value = 42.375801000000003
$ws.Range("H3").Value = -72.519867000000005

# Row 4 - CC-2343 (Florida, USA)
$ws.Range("G4").Value = 27.664826999999999
$ws.Range("H4").Value = -81.515754999999999

# Row 7 - CC-1374 (France)
$ws.Range("G7").Value = 46.227637999999999
$ws.Range("H7").Value = 2.213749

# Row 8 - CC-2932 (North Carolina)
$ws.Range("G8").Value = 35.759574999999998
$ws.Range("H8").Value = -79.019301999999996

# Row 9 - CC-4414 (Colorado)
$ws.Range("G9").Value = 39.550052999999998
$ws.Range("H9").Value = -105.782066

# Row 10 - CC-1373 (South Deerfield, Massachusetts)
$ws.Range("G10").Value = 42.477150000000002
$ws.Range("H10").Value = -72.607950000000002

# Row 11 - CC-1092 (Amherst, Massachusetts, USA)
$ws.Range("G11").Value = 42.375801000000003
$ws.Range("H11").Value = -72.519867000000005

# Row 14 - cc-3980-vhlr-l33: isolation_location corrected from "unknown" to "Amherst, MA"
$ws.Range("E14").Value = "Amherst, MA"
$ws.Range("G14").Value = 42.375801000000003
$ws.Range("H14").Value = -72.519867000000005

# Update view: zoom + new selected cell
$ws.Application.ActiveWindow.Zoom = 109
$ws.Range("J6").Select()
